{"js": "// Commit: \"fix: desarrollo -> trabajo\"\n// The team-name label \"Equipo de desarrollo\" (and one capitalized\n// occurrence \"Equipo de Desarrollo\") is renamed to \"Equipo de trabajo\"\n// everywhere it appears in the document body (table cells throughout the\n// risk-management plan).\n\nconst body = context.document.body;\n\n// --- 1. The single capitalized occurrence: \"Equipo de Desarrollo\" ---\nconst capMatches = body.search(\"Equipo de Desarrollo\", {\n  matchCase: true,\n  matchWholeWord: false\n});\ncapMatches.load(\"items\");\nawait context.sync();\n\nfor (const match of capMatches.items) {\n  match.insertText(\"Equipo de trabajo\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2. The six lower-case occurrences: \"Equipo de desarrollo\" ---\nconst lowerMatches = body.search(\"Equipo de desarrollo\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nlowerMatches.load(\"items\");\nawait context.sync();\n\nfor (const match of lowerMatches.items) {\n  match.insertText(\"Equipo de trabajo\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Commit: \"fix: desarrollo -> trabajo\"\n# The team-name label \"Equipo de desarrollo\" (and one capitalized\n# occurrence \"Equipo de Desarrollo\") is renamed to \"Equipo de trabajo\"\n# everywhere it appears in the document body (table cells throughout the\n# risk-management plan).\n\n$d = $word.ActiveDocument\n\n# --- 1. The single capitalized occurrence: \"Equipo de Desarrollo\" ---\n$findCap = $d.Content.Find\n$findCap.ClearFormatting()\n$findCap.Replacement.ClearFormatting()\n$findCap.Text = \"Equipo de Desarrollo\"\n$findCap.Replacement.Text = \"Equipo de trabajo\"\n$findCap.MatchCase = $true\n$findCap.MatchWholeWord = $false\n$findCap.Execute($findCap.Text, $true, $false, $false, $false, $false, $true, 1, $false, $findCap.Replacement.Text, 2)\n\n# --- 2. The six lower-case occurrences: \"Equipo de desarrollo\" ---\n$findLower = $d.Content.Find\n$findLower.ClearFormatting()\n$findLower.Replacement.ClearFormatting()\n$findLower.Text = \"Equipo de desarrollo\"\n$findLower.Replacement.Text = \"Equipo de trabajo\"\n$findLower.MatchCase = $true\n$findLower.MatchWholeWord = $false\n$findLower.Execute($findLower.Text, $true, $false, $false, $false, $false, $true, 1, $false, $findLower.Replacement.Text, 2)\n\nWrite-Output \"done\"\n"}
